# Update "Pais" (countries) data snapshot + provincias Spain
# This mirrors a periodic refresh of the COVID-19 tracker: the countries
# table (rows sorted descending by "Casos totales") gets new totals for a
# handful of countries - which also re-shuffles their relative order - and
# the "last updated" timestamp advances from 11:22 to 11:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 11:52"

# --- Row 16 ---
$ws.Cells.Item(16, 2).Value = 47334
$ws.Cells.Item(16, 3).Value = 647
$ws.Cells.Item(16, 4).Value = 10943
$ws.Cells.Item(16, 5).Value = 29060
$ws.Cells.Item(16, 6).Value = 876
$ws.Cells.Item(16, 7).Value = 124
$ws.Cells.Item(16, 8).Value = 7331

# --- Row 28 ---
$ws.Cells.Item(28, 2).Value = 15357
$ws.Cells.Item(28, 3).Value = 83
$ws.Cells.Item(28, 4).Value = 12580
$ws.Cells.Item(28, 5).Value = 2208
$ws.Cells.Item(28, 7).Value = 20
$ws.Cells.Item(28, 8).Value = 569

# --- Row 44 ---
$ws.Cells.Item(44, 5).Value = 7361
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 206

# --- Rows 79-84: Eslovenia/Rep. Macedonia/Cuba/Bulgaria/Eslovaquia/Lituania
#     re-sort (by new "Casos totales") and get refreshed case counts ---
$ws.Cells.Item(79, 1).Value = "Eslovenia"
$ws.Cells.Item(79, 2).Value = 1408
$ws.Cells.Item(79, 3).Value = 6
$ws.Cells.Item(79, 4).Value = 223
$ws.Cells.Item(79, 5).Value = 1099
$ws.Cells.Item(79, 6).Value = 24
$ws.Cells.Item(79, 7).Value = 3
$ws.Cells.Item(79, 8).Value = 86

$ws.Cells.Item(80, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(80, 2).Value = 1399
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(80, 4).Value = 553
$ws.Cells.Item(80, 5).Value = 781
$ws.Cells.Item(80, 6).Value = 13
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 65

$ws.Cells.Item(81, 1).Value = "Cuba"
$ws.Cells.Item(81, 2).Value = 1389
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = 525
$ws.Cells.Item(81, 5).Value = 808
$ws.Cells.Item(81, 6).Value = 12
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 56

$ws.Cells.Item(82, 1).Value = "Bulgaria"
$ws.Cells.Item(82, 2).Value = 1387
$ws.Cells.Item(82, 3).Value = 24
$ws.Cells.Item(82, 4).Value = 222
$ws.Cells.Item(82, 5).Value = 1107
$ws.Cells.Item(82, 6).Value = 41
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 58

$ws.Cells.Item(83, 1).Value = "Eslovaquia"
$ws.Cells.Item(83, 2).Value = 1384
$ws.Cells.Item(83, 3).Value = 3
$ws.Cells.Item(83, 4).Value = 423
$ws.Cells.Item(83, 5).Value = 941
$ws.Cells.Item(83, 6).Value = 6
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = 20

$ws.Cells.Item(84, 1).Value = "Lituania"
$ws.Cells.Item(84, 2).Value = 1344
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = 536
$ws.Cells.Item(84, 5).Value = 764
$ws.Cells.Item(84, 6).Value = 17
$ws.Cells.Item(84, 7).Value = 3
$ws.Cells.Item(84, 8).Value = 44

# --- Row 142 (provincias Spain section) ---
$ws.Cells.Item(142, 2).Value = 126
$ws.Cells.Item(142, 3).Value = 2
$ws.Cells.Item(142, 4).Value = 73
